$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data rows 2-37: Fecha (D), Volumen (J), Precio minimo (K), Precio maximo (L), Precio promedio ponderado (M), Precio $/Kg (P)
$ws.Cells.Item(2,4).Value = 44406
$ws.Cells.Item(2,10).Value = 600
$ws.Cells.Item(2,11).Value = 26000
$ws.Cells.Item(2,12).Value = 28000
$ws.Cells.Item(2,13).Value = 27000
$ws.Cells.Item(2,16).Value = 1080
$ws.Cells.Item(3,4).Value = 44455
$ws.Cells.Item(3,10).Value = 800
$ws.Cells.Item(3,11).Value = 28000
$ws.Cells.Item(3,12).Value = 30000
$ws.Cells.Item(3,13).Value = 29000
$ws.Cells.Item(3,16).Value = 1160
$ws.Cells.Item(4,4).Value = 44398
$ws.Cells.Item(4,10).Value = 500
$ws.Cells.Item(4,11).Value = 26000
$ws.Cells.Item(4,12).Value = 28000
$ws.Cells.Item(4,13).Value = 27000
$ws.Cells.Item(4,16).Value = 1080
$ws.Cells.Item(5,4).Value = 44441
$ws.Cells.Item(5,10).Value = 700
$ws.Cells.Item(5,11).Value = 28000
$ws.Cells.Item(5,12).Value = 30000
$ws.Cells.Item(5,13).Value = 29000
$ws.Cells.Item(5,16).Value = 1160
$ws.Cells.Item(6,4).Value = 44420
$ws.Cells.Item(6,10).Value = 700
$ws.Cells.Item(6,11).Value = 27000
$ws.Cells.Item(6,12).Value = 29000
$ws.Cells.Item(6,13).Value = 28000
$ws.Cells.Item(6,16).Value = 1120
$ws.Cells.Item(7,4).Value = 44462
$ws.Cells.Item(7,10).Value = 400
$ws.Cells.Item(7,11).Value = 22000
$ws.Cells.Item(7,12).Value = 23000
$ws.Cells.Item(7,13).Value = 22500
$ws.Cells.Item(7,16).Value = 900
$ws.Cells.Item(8,4).Value = 44371
$ws.Cells.Item(8,10).Value = 500
$ws.Cells.Item(8,11).Value = 28000
$ws.Cells.Item(8,12).Value = 30000
$ws.Cells.Item(8,13).Value = 29000
$ws.Cells.Item(8,16).Value = 1160
$ws.Cells.Item(9,4).Value = 44405
$ws.Cells.Item(9,10).Value = 500
$ws.Cells.Item(9,11).Value = 26000
$ws.Cells.Item(9,12).Value = 28000
$ws.Cells.Item(9,13).Value = 27000
$ws.Cells.Item(9,16).Value = 1080
$ws.Cells.Item(10,4).Value = 44399
$ws.Cells.Item(10,10).Value = 400
$ws.Cells.Item(10,11).Value = 26000
$ws.Cells.Item(10,12).Value = 28000
$ws.Cells.Item(10,13).Value = 27000
$ws.Cells.Item(10,16).Value = 1080
$ws.Cells.Item(11,4).Value = 44364
$ws.Cells.Item(11,10).Value = 200
$ws.Cells.Item(11,11).Value = 28000
$ws.Cells.Item(11,12).Value = 30000
$ws.Cells.Item(11,13).Value = 29000
$ws.Cells.Item(11,16).Value = 1160
$ws.Cells.Item(12,4).Value = 44419
$ws.Cells.Item(12,10).Value = 600
$ws.Cells.Item(12,11).Value = 27000
$ws.Cells.Item(12,12).Value = 29000
$ws.Cells.Item(12,13).Value = 28000
$ws.Cells.Item(12,16).Value = 1120
$ws.Cells.Item(13,4).Value = 44370
$ws.Cells.Item(13,10).Value = 400
$ws.Cells.Item(13,11).Value = 27000
$ws.Cells.Item(13,12).Value = 28000
$ws.Cells.Item(13,13).Value = 27500
$ws.Cells.Item(13,16).Value = 1100
$ws.Cells.Item(14,4).Value = 44412
$ws.Cells.Item(14,10).Value = 600
$ws.Cells.Item(14,11).Value = 25000
$ws.Cells.Item(14,12).Value = 27000
$ws.Cells.Item(14,13).Value = 26000
$ws.Cells.Item(14,16).Value = 1040
$ws.Cells.Item(15,4).Value = 44343
$ws.Cells.Item(15,10).Value = 200
$ws.Cells.Item(15,11).Value = 26000
$ws.Cells.Item(15,12).Value = 28000
$ws.Cells.Item(15,13).Value = 27000
$ws.Cells.Item(15,16).Value = 1080
$ws.Cells.Item(16,4).Value = 44447
$ws.Cells.Item(16,10).Value = 600
$ws.Cells.Item(16,11).Value = 28000
$ws.Cells.Item(16,12).Value = 30000
$ws.Cells.Item(16,13).Value = 29000
$ws.Cells.Item(16,16).Value = 1160
$ws.Cells.Item(17,4).Value = 44356
$ws.Cells.Item(17,10).Value = 300
$ws.Cells.Item(17,11).Value = 26000
$ws.Cells.Item(17,12).Value = 28000
$ws.Cells.Item(17,13).Value = 27000
$ws.Cells.Item(17,16).Value = 1080
$ws.Cells.Item(18,4).Value = 44454
$ws.Cells.Item(18,10).Value = 1000
$ws.Cells.Item(18,11).Value = 28000
$ws.Cells.Item(18,12).Value = 30000
$ws.Cells.Item(18,13).Value = 29000
$ws.Cells.Item(18,16).Value = 1160
$ws.Cells.Item(19,4).Value = 44427
$ws.Cells.Item(19,10).Value = 300
$ws.Cells.Item(19,11).Value = 28000
$ws.Cells.Item(19,12).Value = 30000
$ws.Cells.Item(19,13).Value = 29000
$ws.Cells.Item(19,16).Value = 1160
$ws.Cells.Item(20,4).Value = 44413
$ws.Cells.Item(20,10).Value = 700
$ws.Cells.Item(20,11).Value = 26000
$ws.Cells.Item(20,12).Value = 28000
$ws.Cells.Item(20,13).Value = 27000
$ws.Cells.Item(20,16).Value = 1080
$ws.Cells.Item(21,4).Value = 44426
$ws.Cells.Item(21,10).Value = 400
$ws.Cells.Item(21,11).Value = 28000
$ws.Cells.Item(21,12).Value = 30000
$ws.Cells.Item(21,13).Value = 29000
$ws.Cells.Item(21,16).Value = 1160
$ws.Cells.Item(22,4).Value = 44469
$ws.Cells.Item(22,10).Value = 600
$ws.Cells.Item(22,11).Value = 22000
$ws.Cells.Item(22,12).Value = 24000
$ws.Cells.Item(22,13).Value = 23000
$ws.Cells.Item(22,16).Value = 920
$ws.Cells.Item(23,4).Value = 44434
$ws.Cells.Item(23,10).Value = 500
$ws.Cells.Item(23,11).Value = 28000
$ws.Cells.Item(23,12).Value = 30000
$ws.Cells.Item(23,13).Value = 29000
$ws.Cells.Item(23,16).Value = 1160
$ws.Cells.Item(24,4).Value = 44391
$ws.Cells.Item(24,10).Value = 100
$ws.Cells.Item(24,11).Value = 26000
$ws.Cells.Item(24,12).Value = 28000
$ws.Cells.Item(24,13).Value = 27000
$ws.Cells.Item(24,16).Value = 1080
$ws.Cells.Item(25,4).Value = 44475
$ws.Cells.Item(25,10).Value = 1000
$ws.Cells.Item(25,11).Value = 22000
$ws.Cells.Item(25,12).Value = 24000
$ws.Cells.Item(25,13).Value = 23000
$ws.Cells.Item(25,16).Value = 920
$ws.Cells.Item(26,4).Value = 44349
$ws.Cells.Item(26,10).Value = 600
$ws.Cells.Item(26,11).Value = 26000
$ws.Cells.Item(26,12).Value = 28000
$ws.Cells.Item(26,13).Value = 27000
$ws.Cells.Item(26,16).Value = 1080
$ws.Cells.Item(27,4).Value = 44350
$ws.Cells.Item(27,10).Value = 700
$ws.Cells.Item(27,11).Value = 28000
$ws.Cells.Item(27,12).Value = 30000
$ws.Cells.Item(27,13).Value = 29000
$ws.Cells.Item(27,16).Value = 1160
$ws.Cells.Item(28,4).Value = 44448
$ws.Cells.Item(28,10).Value = 400
$ws.Cells.Item(28,11).Value = 28000
$ws.Cells.Item(28,12).Value = 30000
$ws.Cells.Item(28,13).Value = 29000
$ws.Cells.Item(28,16).Value = 1160
$ws.Cells.Item(29,4).Value = 44385
$ws.Cells.Item(29,10).Value = 500
$ws.Cells.Item(29,11).Value = 26000
$ws.Cells.Item(29,12).Value = 28000
$ws.Cells.Item(29,13).Value = 27000
$ws.Cells.Item(29,16).Value = 1080
$ws.Cells.Item(30,4).Value = 44435
$ws.Cells.Item(30,10).Value = 900
$ws.Cells.Item(30,11).Value = 28000
$ws.Cells.Item(30,12).Value = 30000
$ws.Cells.Item(30,13).Value = 29000
$ws.Cells.Item(30,16).Value = 1160
$ws.Cells.Item(31,4).Value = 44377
$ws.Cells.Item(31,10).Value = 500
$ws.Cells.Item(31,11).Value = 26000
$ws.Cells.Item(31,12).Value = 28000
$ws.Cells.Item(31,13).Value = 27000
$ws.Cells.Item(31,16).Value = 1080
$ws.Cells.Item(32,4).Value = 44433
$ws.Cells.Item(32,10).Value = 400
$ws.Cells.Item(32,11).Value = 28000
$ws.Cells.Item(32,12).Value = 30000
$ws.Cells.Item(32,13).Value = 29000
$ws.Cells.Item(32,16).Value = 1160
$ws.Cells.Item(33,4).Value = 44468
$ws.Cells.Item(33,10).Value = 500
$ws.Cells.Item(33,11).Value = 23000
$ws.Cells.Item(33,12).Value = 25000
$ws.Cells.Item(33,13).Value = 24000
$ws.Cells.Item(33,16).Value = 960
$ws.Cells.Item(34,4).Value = 44384
$ws.Cells.Item(34,10).Value = 400
$ws.Cells.Item(34,11).Value = 26000
$ws.Cells.Item(34,12).Value = 28000
$ws.Cells.Item(34,13).Value = 27000
$ws.Cells.Item(34,16).Value = 1080
$ws.Cells.Item(35,4).Value = 44363
$ws.Cells.Item(35,10).Value = 240
$ws.Cells.Item(35,11).Value = 28000
$ws.Cells.Item(35,12).Value = 30000
$ws.Cells.Item(35,13).Value = 29000
$ws.Cells.Item(35,16).Value = 1160
$ws.Cells.Item(36,4).Value = 44357
$ws.Cells.Item(36,10).Value = 340
$ws.Cells.Item(36,11).Value = 28000
$ws.Cells.Item(36,12).Value = 30000
$ws.Cells.Item(36,13).Value = 29000
$ws.Cells.Item(36,16).Value = 1160
$ws.Cells.Item(37,4).Value = 44392
$ws.Cells.Item(37,10).Value = 100
$ws.Cells.Item(37,11).Value = 26000
$ws.Cells.Item(37,12).Value = 28000
$ws.Cells.Item(37,13).Value = 27000
$ws.Cells.Item(37,16).Value = 1080

# New data row 38 (new weekly record)
$ws.Cells.Item(38,1).Value = 2
$ws.Cells.Item(38,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(38,3).Value = "Coquimbo"
$ws.Cells.Item(38,4).Value = 44461
$ws.Cells.Item(38,4).NumberFormat = $ws.Cells.Item(37,4).NumberFormat
$ws.Cells.Item(38,5).Value = 4
$ws.Cells.Item(38,6).Value = 100112022
$ws.Cells.Item(38,7).Value = "Arveja Verde"
$ws.Cells.Item(38,8).Value = "Perfection"
$ws.Cells.Item(38,9).Value = "Primera"
$ws.Cells.Item(38,10).Value = 500
$ws.Cells.Item(38,11).Value = 23000
$ws.Cells.Item(38,12).Value = 25000
$ws.Cells.Item(38,13).Value = 24000
$ws.Cells.Item(38,14).Value = "$/malla 25 kilos"
$ws.Cells.Item(38,15).Value = "Provincia de Limarí"
$ws.Cells.Item(38,16).Value = 960
$ws.Cells.Item(38,17).Value = 25
$ws.Cells.Item(38,18).Value = "Hortaliza"
